$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "routes": insert a new row for "/logout" right after the "/login"
# row, mark login as done/tested, and re-point the remarks accordingly.
# ---------------------------------------------------------------------------
$routes = $wb.Worksheets.Item("routes")

# Insert a fresh row 3 (pushes /register and everything below down by one).
$routes.Rows.Item(3).Insert() | Out-Null

# Update the existing /login row (row 2): it is now considered done/tested.
$routes.Range("F2").Value = "done"
$routes.Range("G2").Value = "9/23/2024"
$routes.Range("H2").Value = "login tested for all users"

# Fill in the new /logout row (row 3).
$routes.Range("B3").Value = "/logout"
$routes.Range("D3").Value = "logout user from the system"
$routes.Range("C3").Value = "Self"
$routes.Range("E3").Value = "Jasdeep"
$routes.Range("F3").Value = "current task"
$routes.Range("G3").Value = "9/23/2024"

# ---------------------------------------------------------------------------
# Sheet "tasks": log the new login testing task performed on 9/23/2024.
# ---------------------------------------------------------------------------
$tasks = $wb.Worksheets.Item("tasks")

# Match the date formatting already used by the column (copy A6's format).
$tasks.Range("A6").Copy()
$tasks.Range("A7").PasteSpecial(-4122) | Out-Null
$tasks.Range("A7").Value = "9/23/2024"
$tasks.Range("B7").Value = "Jasdeep"
$tasks.Range("C7").Value = "login tested with new User model"
$tasks.Range("D7").Value = "Tested with jest and supertest. api doc to be created."

# Back to routes to finish the remarks column, then leave the selection
# where the author's session ended up.
$routes.Range("H3").Value = "to check best technique to use with JWT "

$tasks.Range("D7").Select() | Out-Null
$routes.Range("H4").Select() | Out-Null
